$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.987.89"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "2.325.76"
$ws.Range("E3").Value = "  -3.87%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.103"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.01%  "
$ws.Range("D14").Value = "2.753.39"
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "59.029.17"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "2.356.79"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "317.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.168"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.26%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.391"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "307.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "143.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0951"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0499"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.561"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0213"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.939"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "
